$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = "Thai Red Curry Fried Rice"

# Row 3
$ws.Range("A3").Value = "Lasagne"
$ws.Range("C3").Value = "Veggie Chow Mein"

# Row 4
$ws.Range("A4").Value = "Feta Chicken Pasta"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "Chicken and Vegetable Nimono"

# Row 5
$ws.Range("A5").Value = "Spicy Vegetable Curry"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Sticky Asian Meatballs"

# Row 6
$ws.Range("A6").Value = "Citrus Chicken and Vegetables"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "General Tso's Tofu"

# Row 7
$ws.Range("A7").Value = "Veggie Lasagne"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "Tonkatsu Pork"

# Row 8
$ws.Range("A8").Value = "Bacon Stuffed Mushrooms"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "White Bean Chicken"

# Row 9
$ws.Range("A9").Value = "Spring Vegetable Risotto"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "Plum Chicken"

# Row 10
$ws.Range("A10").Value = "Beef Stroganoff"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "Beef Satay Skewers"

# Row 11
$ws.Range("A11").Value = "Pesto Penne Pasta"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = "Braised Pork Mince"

# Row 15
$ws.Range("A15").Value = "Tarragon Chicken"

# Row 16
$ws.Range("A16").Value = "Tonkatsu Pork"

# Row 17 - becomes numeric 0/0
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 0

# Row 18
$ws.Range("A18").Value = "Tonkatsu Pork"
$ws.Range("B18").Value = 1

# Row 19
$ws.Range("A19").Value = "Thai Red Curry Fried Rice"
$ws.Range("B19").Value = 1

# Row 20
$ws.Range("A20").Value = "Beef Burrito Bowl"

# Row 21 - becomes numeric 0/0
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = 0

# Row 22 - becomes numeric 0/0
$ws.Range("A22").Value = 0
$ws.Range("B22").Value = 0

# Row 23 - becomes numeric 0/0
$ws.Range("A23").Value = 0
$ws.Range("B23").Value = 0
